$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text (preserving trailing zeros / formatting).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '43.841.51'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '2.187.46'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '293.52'
$ws.Range("E5").Value = '  -4.06%  '
$ws.Range("D6").Value = '87.29'
$ws.Range("E6").Value = '  -6.26%  '
$ws.Range("D7").Value = '0.562'
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("D9").Value = '0.477'
$ws.Range("E9").Value = '  -8.66%  '
$ws.Range("D10").Value = '31.96'
$ws.Range("E10").Value = '  -7.16%  '
$ws.Range("D11").Value = '0.0760'
$ws.Range("E11").Value = '  -6.11%  '
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '6.70'
$ws.Range("E13").Value = '  -6.02%  '
$ws.Range("D14").Value = '2.516.50'
$ws.Range("E14").Value = '  -2.68%  '
$ws.Range("D15").Value = '2.256.30'
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").Value = '12.85'
$ws.Range("E16").Value = '  -5.24%  '
$ws.Range("D17").Value = '0.759'
$ws.Range("E17").Value = '  -9.19%  '
$ws.Range("D18").Value = '43.297.67'
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '0.0₃0872'
$ws.Range("E19").Value = '  -9.20%  '
$ws.Range("D20").Value = '5.78'
$ws.Range("E20").Value = '  -8.88%  '
$ws.Range("D21").Value = '10.63'
$ws.Range("E21").Value = '  -13.74%  '
$ws.Range("D22").Value = '62.29'
$ws.Range("E22").Value = '  -5.08%  '
$ws.Range("D23").Value = '227.59'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("D24").Value = '2.75'
$ws.Range("E24").Value = '  -6.36%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '1.80'
$ws.Range("E26").Value = '  -8.67%  '
$ws.Range("D27").Value = '2.22'
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '35.27'
$ws.Range("E28").Value = '  -8.81%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  -7.19%  '
$ws.Range("D30").Value = '18.97'
$ws.Range("E30").Value = '  -5.28%  '
$ws.Range("D31").Value = '145.56'
$ws.Range("E31").Value = '  -5.27%  '
$ws.Range("D32").Value = '5.23'
$ws.Range("E32").Value = '  -11.62%  '
$ws.Range("D33").Value = '2.49'
$ws.Range("E33").Value = '  -6.09%  '
$ws.Range("D34").Value = '0.0721'
$ws.Range("E34").Value = '  -9.49%  '
$ws.Range("D35").Value = '0.115'
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").Value = '2.87'
$ws.Range("E36").Value = '  -7.74%  '
$ws.Range("D37").Value = '0.100'
$ws.Range("E37").Value = '  -7.31%  '
$ws.Range("D38").Value = '1.63'
$ws.Range("E38").Value = '  -10.54%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0277'
$ws.Range("E39").Value = '  -7.76%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '3.46'
$ws.Range("E40").Value = '  -9.10%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Value = '13.09'
$ws.Range("E42").Value = '  -9.95%  '
$ws.Range("D43").Value = '3.01'
$ws.Range("E43").Value = '  -12.33%  '
$ws.Range("D44").Value = '1.749.18'
$ws.Range("E44").Value = '  +0.96%  '
$ws.Range("D45").Value = '1.63'
$ws.Range("E45").Value = '  +2.78%  '
$ws.Range("D46").Value = '72.12'
$ws.Range("E46").Value = '  -10.09%  '
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").Value = '65.36'
$ws.Range("E47").Value = '  -6.05%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.170'
$ws.Range("E48").Value = '  -11.89%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '90.65'
$ws.Range("E49").Value = '  -8.62%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.399.16'
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.70'
$ws.Range("E51").Value = '  +7.23%  '
